$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# (values like "1.007" or "102.00" must not be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.928.22"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "1.846.33"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "309.88"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4733"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "0.07224"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "0.9237"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").Value = "19.63"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "0.07624"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "1.882.72"
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("D14").Value = "5.313"
$ws.Range("D16").Value = "88.34"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "0.000008671"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "26.956.31"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "14.57"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").Value = "5.038"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "10.67"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "1.916"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "152.27"
$ws.Range("D26").Value = "18.17"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "2.002"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "114.36"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "4.991"
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("D30").Value = "0.08842"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "3.283"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").Value = "0.7475"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("D33").Value = "1.167"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "2.767"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.502"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "0.01951"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "0.05259"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("D39").Value = "2.969"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "0.5206"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").Value = "6.917"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "0.1511"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "8.205"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "10.56"
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").Value = "0.4698"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").Value = "1.009"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "102.00"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").Value = "1.601"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").Value = "65.49"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "0.06035"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "0.8853"
$ws.Range("E51").Value = "  +4.18%  "
